{"js": "// Mark the \"\u0418\u043c\u043f\u043e\u0440\u0442 \u0442\u043e\u0432\u0430\u0440\u043e\u0432 \u0438\u0437 \u0444\u0430\u0439\u043b\u0430\" task as done by applying\n// strikethrough formatting to the whole paragraph (paragraph mark\n// included), matching the other completed tasks in the list.\n\nconst body = context.document.body;\n\n// Locate the paragraph by its text so the edit is anchored to the\n// content described in the change, not a fragile numeric index.\nconst results = body.search(\"\u0418\u043c\u043f\u043e\u0440\u0442 \u0442\u043e\u0432\u0430\u0440\u043e\u0432 \u0438\u0437 \u0444\u0430\u0439\u043b\u0430\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Paragraph \"\u0418\u043c\u043f\u043e\u0440\u0442 \u0442\u043e\u0432\u0430\u0440\u043e\u0432 \u0438\u0437 \u0444\u0430\u0439\u043b\u0430\" not found.');\n}\n\nconst targetRange = results.items[0];\nconst paragraph = targetRange.paragraphs.getFirst();\n\n// Setting font.strikeThrough on the paragraph applies <w:strike/> to\n// both the run(s) of text and the paragraph mark's run properties\n// (w:pPr/w:rPr), exactly like the already-completed sibling tasks.\nparagraph.font.strikeThrough = true;\n\nawait context.sync();\n", "ps1": "# Mark the \"\u0418\u043c\u043f\u043e\u0440\u0442 \u0442\u043e\u0432\u0430\u0440\u043e\u0432 \u0438\u0437 \u0444\u0430\u0439\u043b\u0430\" task as done by applying\n# strikethrough formatting to the whole paragraph (paragraph mark\n# included), matching the other completed tasks in the list.\n\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$range.Find.ClearFormatting()\n$found = $range.Find.Execute([ref]\"\u0418\u043c\u043f\u043e\u0440\u0442 \u0442\u043e\u0432\u0430\u0440\u043e\u0432 \u0438\u0437 \u0444\u0430\u0439\u043b\u0430\")\n\nif ($found) {\n    # Expand to the full paragraph (including the paragraph mark) so\n    # both the run and the pPr/rPr of the paragraph get w:strike,\n    # exactly like the sibling \"done\" tasks elsewhere in the list.\n    $para = $range.Paragraphs(1).Range\n    $para.Font.StrikeThrough = 1\n}\n"}
